$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.859.57"
$ws.Range("E2").Value = "  -1.25%  "
$ws.Range("D3").Value = "'1.903.57"
$ws.Range("E3").Value = "  -0.64%  "
$ws.Range("E4").Value = "  -0.59%  "
$ws.Range("D5").Value = "'313.02"
$ws.Range("E5").Value = "  -1.35%  "
$ws.Range("D6").Value = "'0.9999"
$ws.Range("E6").Value = "  -0.53%  "
$ws.Range("D7").Value = "'0.4946"
$ws.Range("E7").Value = "  +2.16%  "
$ws.Range("D8").Value = "'0.3821"
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").Value = "'0.07334"
$ws.Range("E9").Value = "  -0.66%  "
$ws.Range("D10").Value = "'0.9099"
$ws.Range("E10").Value = "  -3.24%  "
$ws.Range("D11").Value = "'20.98"
$ws.Range("E11").Value = "  +0.52%  "
$ws.Range("D12").Value = "'0.07621"
$ws.Range("E12").Value = "  -2.64%  "
$ws.Range("D13").Value = "'1.913.84"
$ws.Range("E13").Value = "  +0.06%  "
$ws.Range("D14").Value = "'5.482"
$ws.Range("E14").Value = "  -0.43%  "
$ws.Range("D15").Value = "'6.639"
$ws.Range("E15").Value = "  -0.13%  "
$ws.Range("D16").Value = "'91.35"
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("E17").Value = "  -0.55%  "
$ws.Range("D18").Value = "'0.000008732"
$ws.Range("E18").Value = "  -1.24%  "
$ws.Range("D19").Value = "'0.9995"
$ws.Range("E19").Value = "  -0.51%  "
$ws.Range("D20").Value = "'27.895.18"
$ws.Range("E20").Value = "  -1.19%  "
$ws.Range("D22").Value = "'5.134"
$ws.Range("E22").Value = "  -0.55%  "
$ws.Range("D23").Value = "'10.80"
$ws.Range("E23").Value = "  -1.17%  "
$ws.Range("D24").Value = "'154.65"
$ws.Range("E24").Value = "  -1.09%  "
$ws.Range("D25").Value = "'1.867"
$ws.Range("E25").Value = "  -3.01%  "
$ws.Range("D26").Value = "'2.227"
$ws.Range("E26").Value = "  +5.98%  "
$ws.Range("D27").Value = "'18.42"
$ws.Range("E27").Value = "  -0.88%  "
$ws.Range("D28").Value = "'115.28"
$ws.Range("E28").Value = "  -0.96%  "
$ws.Range("D29").Value = "'4.927"
$ws.Range("E29").Value = "  -0.90%  "
$ws.Range("D30").Value = "'0.08948"
$ws.Range("E30").Value = "  +0.49%  "
$ws.Range("D31").Value = "'3.211"
$ws.Range("E31").Value = "  -4.34%  "
$ws.Range("D32").Value = "'1.240"
$ws.Range("E32").Value = "  -0.94%  "
$ws.Range("E33").Value = "  -0.27%  "
$ws.Range("D34").Value = "'4.644"
$ws.Range("E34").Value = "  -1.29%  "
$ws.Range("D35").Value = "'0.02064"
$ws.Range("E35").Value = "  +0.52%  "
$ws.Range("D36").Value = "'2.576"
$ws.Range("E36").Value = "  -3.38%  "
$ws.Range("D37").Value = "'1.099"
$ws.Range("E37").Value = "  -0.26%  "
$ws.Range("E38").Value = "  +0.07%  "
$ws.Range("D39").Value = "'0.05294"
$ws.Range("E39").Value = "  -0.53%  "
$ws.Range("D40").Value = "'3.006"
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("D41").Value = "'6.994"
$ws.Range("E41").Value = "  -0.85%  "
$ws.Range("D42").Value = "'8.554"
$ws.Range("E42").Value = "  +1.03%  "
$ws.Range("D43").Value = "'0.1522"
$ws.Range("E43").Value = "  -0.54%  "
$ws.Range("D44").Value = "'111.16"
$ws.Range("E44").Value = "  +3.91%  "
$ws.Range("D45").Value = "'10.65"
$ws.Range("E45").Value = "  -0.55%  "
$ws.Range("D46").Value = "'0.4797"
$ws.Range("E46").Value = "  -1.36%  "
$ws.Range("D47").Value = "'0.9996"
$ws.Range("E47").Value = "  -0.56%  "
$ws.Range("D48").Value = "'1.641"
$ws.Range("E48").Value = "  -1.15%  "
$ws.Range("D49").Value = "'67.54"
$ws.Range("E49").Value = "  -1.59%  "
$ws.Range("D50").Value = "'0.06074"
$ws.Range("E50").Value = "  -0.67%  "
$ws.Range("D51").Value = "'0.8992"
$ws.Range("E51").Value = "  -1.16%  "
